$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)

# Sheet 1
$ws1.Range("F2").Value = 609
$ws1.Range("F4").Value = 6456
$ws1.Range("F5").Value = 737
$ws1.Range("F7").Value = 82
$ws1.Range("F8").Value = 543
$ws1.Range("F11").Value = 721
$ws1.Range("F12").Value = 1202
$ws1.Range("F14").Value = 87
$ws1.Range("F16").Value = 448
$ws1.Range("F20").Value = 677
$ws1.Range("F21").Value = 394
$ws1.Range("F22").Value = 404
$ws1.Range("F25").Value = 164
$ws1.Range("F26").Value = 2230
$ws1.Range("F29").Value = 404
$ws1.Range("F31").Value = 3626
$ws1.Range("F33").Value = 642

# Sheet 2
$ws2.Range("F12").Value = 1022
$ws2.Range("F14").Value = 112
$ws2.Range("F17").Value = 4
$ws2.Range("F25").Value = 19
$ws2.Range("F27").Value = 198
$ws2.Range("F31").Value = 213
$ws2.Range("F35").Value = 1670
$ws2.Range("F37").Value = 6

# Sheet 3
$ws3.Range("F4").Value = 1204
$ws3.Range("F6").Value = 1584
$ws3.Range("F10").Value = 799

# Sheet 4
$ws4.Range("F4").Value = 1204
$ws4.Range("F5").Value = 1584
$ws4.Range("F8").Value = 799
$ws4.Range("F9").Value = 609
$ws4.Range("F10").Value = 6456
$ws4.Range("F12").Value = 737
$ws4.Range("F15").Value = 82
$ws4.Range("F16").Value = 543
$ws4.Range("F19").Value = 721
$ws4.Range("F21").Value = 112
$ws4.Range("F22").Value = 112
$ws4.Range("F24").Value = 1202
$ws4.Range("F25").Value = 87
$ws4.Range("F31").Value = 19
$ws4.Range("F32").Value = 677
$ws4.Range("F33").Value = 394
$ws4.Range("F34").Value = 404
$ws4.Range("F38").Value = 164
$ws4.Range("F39").Value = 2230
$ws4.Range("F40").Value = 213
$ws4.Range("F43").Value = 1670
$ws4.Range("F44").Value = 1670
$ws4.Range("F46").Value = 404
$ws4.Range("F47").Value = 3626
$ws4.Range("F48").Value = 6
$ws4.Range("F51").Value = 642
